$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H70").Value = 3003.125
$ws.Range("I70").Value = 997.5
$ws.Range("J70").Value = 3671.6667
$ws.Range("K70").Value = 2992.5
$ws.Range("L70").Value = 11015.0001
$ws.Range("M70").Value = -2722.5
$ws.Range("N70").Value = -11555.0001
$ws.Range("H73").Value = 3003.125
$ws.Range("I73").Value = 997.5
$ws.Range("J73").Value = 3671.6667
$ws.Range("K73").Value = 2992.5
$ws.Range("L73").Value = 11015.0001
$ws.Range("M73").Value = -2056.5
$ws.Range("N73").Value = -12887.0001
$ws.Range("H127").Value = 1028
$ws.Range("I127").Value = 1028
$ws.Range("K127").Value = 3084
$ws.Range("M127").Value = 1876

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 39475.652
$ws.Range("J52").Value = 39475.652
$ws.Range("L52").Value = 39475.652
$ws.Range("N52").Value = -40111.652
$ws.Range("H74").Value = 929.4545000000001
$ws.Range("I74").Value = 813.75
$ws.Range("J74").Value = 1238
$ws.Range("K74").Value = 813.75
$ws.Range("L74").Value = 1238
$ws.Range("M74").Value = 60.25
$ws.Range("N74").Value = -2986
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("N75").ClearContents()
$ws.Range("H77").Value = 929.4545000000001
$ws.Range("I77").Value = 813.75
$ws.Range("J77").Value = 1238
$ws.Range("K77").Value = 4068.75
$ws.Range("L77").Value = 6190
$ws.Range("M77").Value = 299.25
$ws.Range("N77").Value = -14926
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("N78").ClearContents()
$ws.Range("H110").Value = 1594.5
$ws.Range("I110").Value = 717.1818
$ws.Range("K110").Value = 717.1818
$ws.Range("M110").Value = 1327.8182
$ws.Range("H122").Value = 2426.1667
$ws.Range("J122").Value = 3353.0833
$ws.Range("L122").Value = 10059.2499
$ws.Range("N122").Value = -14959.2499
$ws.Range("H132").Value = 35720068
$ws.Range("I132").Value = 62506340
$ws.Range("J132").Value = 5035.5
$ws.Range("K132").Value = 187519020
$ws.Range("L132").Value = 15106.5
$ws.Range("M132").Value = -187516490
$ws.Range("N132").Value = -20166.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4450.4546
$ws.Range("I62").Value = 2959.1667
$ws.Range("J62").Value = 6240
$ws.Range("K62").Value = 2959.1667
$ws.Range("L62").Value = 6240
$ws.Range("M62").Value = -2335.1667
$ws.Range("N62").Value = -7488
$ws.Range("H65").Value = 4450.4546
$ws.Range("I65").Value = 2959.1667
$ws.Range("J65").Value = 6240
$ws.Range("K65").Value = 14795.8335
$ws.Range("L65").Value = 31200
$ws.Range("M65").Value = -11675.8335
$ws.Range("N65").Value = -37440
$ws.Range("H122").Value = 1890.4138
$ws.Range("I122").Value = 1669.8948
$ws.Range("J122").Value = 2309.4
$ws.Range("K122").Value = 5009.6844
$ws.Range("L122").Value = 6928.200000000001
$ws.Range("M122").Value = -2559.6844
$ws.Range("N122").Value = -11828.2
$ws.Range("H132").Value = 3701.9
$ws.Range("I132").Value = 2553.85
$ws.Range("J132").Value = 5998
$ws.Range("K132").Value = 7661.549999999999
$ws.Range("L132").Value = 17994
$ws.Range("M132").Value = -5131.549999999999
$ws.Range("N132").Value = -23054

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 52059.035
$ws.Range("I131").Value = 4633.3335
$ws.Range("J131").Value = 54795.133
$ws.Range("K131").Value = 13900.0005
$ws.Range("L131").Value = 164385.399
$ws.Range("M131").Value = -8860.000499999998
$ws.Range("N131").Value = -174465.399
$ws.Range("H136").Value = 2018.4
$ws.Range("I136").Value = 1558.6875
$ws.Range("J136").Value = 3857.25
$ws.Range("K136").Value = 4676.0625
$ws.Range("L136").Value = 11571.75
$ws.Range("M136").Value = 423.9375
$ws.Range("N136").Value = -21771.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2518.2334
$ws.Range("I102").Value = 1897.6957
$ws.Range("J102").Value = 4557.143
$ws.Range("K102").Value = 1897.6957
$ws.Range("L102").Value = 4557.143
$ws.Range("M102").Value = -275.6957
$ws.Range("N102").Value = -7801.143
$ws.Range("H126").Value = 3828.1428
$ws.Range("I126").Value = 1864.6666
$ws.Range("J126").Value = 4363.636
$ws.Range("K126").Value = 5593.9998
$ws.Range("L126").Value = 13090.908
$ws.Range("M126").Value = -3123.9998
$ws.Range("N126").Value = -18030.908
$ws.Range("H132").Value = 3487.5334
$ws.Range("I132").Value = 3129.0715
$ws.Range("J132").Value = 3801.1875
$ws.Range("K132").Value = 9387.2145
$ws.Range("L132").Value = 11403.5625
$ws.Range("M132").Value = -6857.2145
$ws.Range("N132").Value = -16463.5625

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 125001550
$ws.Range("I22").Value = 1000000000
$ws.Range("J22").Value = 1778
$ws.Range("K22").Value = 1000000000
$ws.Range("L22").Value = 1778
$ws.Range("M22").Value = -999999705
$ws.Range("N22").Value = -2368
$ws.Range("H27").Value = 125001550
$ws.Range("I27").Value = 1000000000
$ws.Range("J27").Value = 1778
$ws.Range("K27").Value = 1000000000
$ws.Range("L27").Value = 1778
$ws.Range("M27").Value = -999999893
$ws.Range("N27").Value = -1992
$ws.Range("H40").Value = 2413.7
$ws.Range("I40").Value = 1000
$ws.Range("J40").Value = 2767.125
$ws.Range("K40").Value = 1000
$ws.Range("L40").Value = 2767.125
$ws.Range("M40").Value = -864
$ws.Range("N40").Value = -3039.125
$ws.Range("H132").Value = 3113.3225
$ws.Range("I132").Value = 1744.375
$ws.Range("K132").Value = 5233.125
$ws.Range("M132").Value = -2703.125

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 34000
$ws.Range("J80").Value = 34000
$ws.Range("L80").Value = 34000
$ws.Range("N80").Value = -35996
$ws.Range("H83").Value = 34000
$ws.Range("J83").Value = 34000
$ws.Range("L83").Value = 102000
$ws.Range("N83").Value = -111984
$ws.Range("H126").Value = 8336443
$ws.Range("I126").Value = 2413.5
$ws.Range("J126").Value = 25004502
$ws.Range("K126").Value = 7240.5
$ws.Range("L126").Value = 75013506
$ws.Range("M126").Value = -4770.5
$ws.Range("N126").Value = -75018446
